# The underlying source data feeding this export was re-sorted upstream.
# Net visible effect in the sheet: several row-pairs each have their whole
# record (every column A:AY) swapped with their partner row, while the row
# number itself stays fixed. Pairs affected: (5,6), (19,20), (22,23), (24,25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = "AY"
$pairs = @(
    @(5, 6),
    @(19, 20),
    @(22, 23),
    @(24, 25)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("A" + $r1 + ":" + $lastCol + $r1)
    $range2 = $ws.Range("A" + $r2 + ":" + $lastCol + $r2)

    # A few columns hold digit-only strings stored as *text* in the source
    # (Antal "I", and the yyyy-mm-dd "Startdatum"/"Slutdatum" Y/AA). Reading
    # them through .Value/.Value2 as part of the bulk row swap makes Excel
    # re-type them as real numbers/dates, so capture their raw text first
    # and restore it verbatim (forcing Text format) after the swap.
    $textCols = @("I", "Y", "AA")
    $orig1 = @{}
    $orig2 = @{}
    foreach ($col in $textCols) {
        $orig1[$col] = $ws.Range($col + $r1).Value2()
        $orig2[$col] = $ws.Range($col + $r2).Value2()
    }

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1

    foreach ($col in $textCols) {
        $ws.Range($col + $r1).NumberFormat = "@"
        $ws.Range($col + $r2).NumberFormat = "@"
        $ws.Range($col + $r1).Value = $orig2[$col]
        $ws.Range($col + $r2).Value = $orig1[$col]
    }
}
